# Update the "Edges" sheet to make the graph bidirectional: for every
# existing directed edge (A,B) add the reverse edge (B,A), and double the
# edge-count cell A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Edges" sheet is the active/tabSelected sheet

# Update edge count in A1 (11 -> 22)
$ws.Range("A1").Value = 22

# Append reversed edges for rows 2-12 (original edges) as rows 13-23
$edges = @(
    @(2,1),
    @(3,1),
    @(4,1),
    @(4,2),
    @(5,2),
    @(6,3),
    @(8,3),
    @(7,4),
    @(8,5),
    @(6,7),
    @(6,8)
)

$row = 13
foreach ($edge in $edges) {
    $ws.Cells.Item($row, 1).Value = $edge[0]
    $ws.Cells.Item($row, 2).Value = $edge[1]
    $row++
}

# Move the active selection to A2, matching the post-edit workbook state.
$ws.Range("A2").Select()
